$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3864.003151314854
$ws.Range("C3").Value = 4032.565189692972
$ws.Range("C4").Value = 4172.797737398093
$ws.Range("C5").Value = 4275.681424680156
$ws.Range("C6").Value = 4421.103665140147
$ws.Range("C7").Value = 4562.8138225845
$ws.Range("C8").Value = 4168.75072638867
$ws.Range("C9").Value = 4296.681196381269
$ws.Range("C10").Value = 4399.909968231253
$ws.Range("C11").Value = 4339.081639163523
$ws.Range("C12").Value = 4532.3797990523
$ws.Range("C13").Value = 4551.220592847516
$ws.Range("C14").Value = 4632.024569265814
$ws.Range("C15").Value = 4726.52469657353
$ws.Range("C16").Value = 4801.830961006282
$ws.Range("C17").Value = 4823.803935421367
$ws.Range("C18").Value = 4899.702323349131
$ws.Range("C19").Value = 4949.101679397882
$ws.Range("C20").Value = 4985.398429293686
$ws.Range("C21").Value = 5039.45135531525
$ws.Range("C22").Value = 5115.05939331125
$ws.Range("C23").Value = 5194.536377584491
$ws.Range("C24").Value = 5267.091621051989
$ws.Range("C25").Value = 5363.216622853107
$ws.Range("C26").Value = 5376.028799266796
$ws.Range("C27").Value = 5431.134417735461
$ws.Range("C28").Value = 5484.794636495427
$ws.Range("C29").Value = 5532.381084508597
$ws.Range("C30").Value = 5568.857852270792
$ws.Range("C31").Value = 5608.239449197749
$ws.Range("C32").Value = 5656.319138321764
$ws.Range("C33").Value = 5694.667115556107
$ws.Range("C34").Value = 5731.382812788844
$ws.Range("C35").Value = 5732.444684351412
$ws.Range("C36").Value = 5738.771914396423
$ws.Range("C37").Value = 5779.697968151032
$ws.Range("C38").Value = 5770.874073456223
$ws.Range("C39").Value = 5805.17458340013
$ws.Range("C40").Value = 5826.897384379914
$ws.Range("C41").Value = 5891.684674193875
$ws.Range("C42").Value = 5871.974135824853
$ws.Range("C43").Value = 5889.944781996448
$ws.Range("C44").Value = 5913.282604318892
$ws.Range("C45").Value = 5919.938361341304
$ws.Range("C46").Value = 5939.355236882049
$ws.Range("C47").Value = 5988.324192142125
$ws.Range("C48").Value = 5981.220819227029
$ws.Range("C49").Value = 6002.397726965758
$ws.Range("C50").Value = 6025.023720613532
$ws.Range("C51").Value = 6058.582142181468
$ws.Range("C52").Value = 6078.379926516941
$ws.Range("C53").Value = 6113.784192456451
$ws.Range("C54").Value = 6142.145992182102
$ws.Range("C55").Value = 6177.3733886481
$ws.Range("C56").Value = 6216.909589940034
$ws.Range("C57").Value = 6213.433865934981
$ws.Range("C58").Value = 6247.879002986622
$ws.Range("C59").Value = 6276.06097900224
